# Update "想去人数" (want-to-go count) figures in column F across the
# 展览 (sheet1), 演出 (sheet2) and 全部类型 (sheet4) worksheets to match
# the refreshed data output.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1253
$ws1.Range("F5").Value = 923
$ws1.Range("F6").Value = 1669
$ws1.Range("F8").Value = 1121
$ws1.Range("F9").Value = 44
$ws1.Range("F11").Value = 245
$ws1.Range("F12").Value = 11
$ws1.Range("F14").Value = 604
$ws1.Range("F15").Value = 116
$ws1.Range("F20").Value = 68
$ws1.Range("F21").Value = 628
$ws1.Range("F22").Value = 617
$ws1.Range("F23").Value = 105
$ws1.Range("F25").Value = 825
$ws1.Range("F27").Value = 44
$ws1.Range("F29").Value = 234
$ws1.Range("F31").Value = 5

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 302
$ws2.Range("F5").Value = 7
$ws2.Range("F7").Value = 227
$ws2.Range("F11").Value = 17

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1253
$ws4.Range("F6").Value = 923
$ws4.Range("F7").Value = 1669
$ws4.Range("F9").Value = 1121
$ws4.Range("F10").Value = 44
$ws4.Range("F13").Value = 245
$ws4.Range("F14").Value = 11
$ws4.Range("F16").Value = 604
$ws4.Range("F17").Value = 116
$ws4.Range("F21").Value = 302
$ws4.Range("F23").Value = 7
$ws4.Range("F26").Value = 227
$ws4.Range("F27").Value = 227
$ws4.Range("F28").Value = 68
$ws4.Range("F29").Value = 628
$ws4.Range("F30").Value = 617
$ws4.Range("F31").Value = 105
$ws4.Range("F33").Value = 825
$ws4.Range("F36").Value = 44
$ws4.Range("F38").Value = 234
$ws4.Range("F43").Value = 5
$ws4.Range("F44").Value = 17
